# Insert a new data row at row 10 (pushes the existing rows 10..89 down to 11..90,
# matching the price-history sheet growing by one new weekly observation).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(10).Insert()

$fecha = Get-Date -Year 2022 -Month 5 -Day 4 -Hour 0 -Minute 0 -Second 0

$ws.Range("A10").Value = 2
$ws.Range("B10").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C10").Value = "Coquimbo"
$ws.Range("D10").Value = $fecha
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 100112030
$ws.Range("G10").Value = "Poroto granado"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 600
$ws.Range("K10").Value = 16000
$ws.Range("L10").Value = 18000
$ws.Range("M10").Value = 17000
$ws.Range("N10").Value = '$/malla 25 kilos'
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 680
$ws.Range("Q10").Value = 25
$ws.Range("R10").Value = "Hortaliza"
